$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos - updated daily figures ---
$ws.Range("B4").Value = 651310
$ws.Range("C4").Value = 3162
$ws.Range("D4").Value = 56489
$ws.Range("E4").Value = 561489
$ws.Range("G4").Value = 744
$ws.Range("H4").Value = 33332

# --- Row 24: updated figures ---
$ws.Range("D24").Value = 550
$ws.Range("E24").Value = 10657

# --- Rows 29/30: Polonia overtakes Ecuador in the ranking ---
# Row 29 becomes Polonia with its new (higher) totals
$ws.Range("A29").Value = "Polonia"
$ws.Range("B29").Value = 7918
$ws.Range("C29").Value = 336
$ws.Range("D29").Value = 774
$ws.Range("E29").Value = 6830
$ws.Range("F29").Value = 160
$ws.Range("G29").Value = 28
$ws.Range("H29").Value = 314

# Row 30 becomes Ecuador, keeping its previous totals
$ws.Range("A30").Value = "Ecuador"
$ws.Range("B30").Value = 7858
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 780
$ws.Range("E30").Value = 6690
$ws.Range("F30").Value = 135
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 388

# --- Row 33: updated figure ---
$ws.Range("F33").Value = 92

# --- Row 37: updated figures ---
$ws.Range("B37").Value = 6359
$ws.Range("C37").Value = 58
$ws.Range("D37").Value = 972
$ws.Range("E37").Value = 5218
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 169

# --- Row 70: updated figures ---
$ws.Range("D70").Value = 277
$ws.Range("E70").Value = 1068

# --- Row 81: updated figure ---
$ws.Range("F81").Value = 16
